$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.041.56"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "1.762.25"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.49"
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3927"
$ws.Range("E7").Value = "  +2.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3396"
$ws.Range("E8").Value = "  -1.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.40"
$ws.Range("E9").Value = "  -3.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.122"
$ws.Range("E10").Value = "  -2.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07234"
$ws.Range("E11").Value = "  -2.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.41"
$ws.Range("E13").Value = "  -3.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.152"
$ws.Range("E14").Value = "  -4.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.110"
$ws.Range("E15").Value = "  -3.46%  "
$ws.Range("D16").Value = "1.733.60"
$ws.Range("E16").Value = "  -3.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001061"
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06620"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.45"
$ws.Range("E19").Value = "  -2.21%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.99"
$ws.Range("E21").Value = "  -2.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.235"
$ws.Range("E22").Value = "  -3.40%  "
$ws.Range("D23").Value = "28.025.06"
$ws.Range("E23").Value = "  -0.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.68"
$ws.Range("E24").Value = "  -3.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.398"
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.51"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.96"
$ws.Range("E27").Value = "  -3.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.321"
$ws.Range("E28").Value = "  -4.42%  "
$ws.Range("D29").Value = "1.956.83"
$ws.Range("E29").Value = "  -1.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.277"
$ws.Range("E30").Value = "  -11.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "129.54"
$ws.Range("E31").Value = "  -4.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.080"
$ws.Range("E32").Value = "  +1.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.825"
$ws.Range("E33").Value = "  -4.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08737"
$ws.Range("E34").Value = "  -2.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.10"
$ws.Range("E35").Value = "  -5.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06194"
$ws.Range("E36").Value = "  -3.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02292"
$ws.Range("E37").Value = "  -5.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.150"
$ws.Range("E38").Value = "  -4.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6509"
$ws.Range("E39").Value = "  -5.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2116"
$ws.Range("E40").Value = "  -2.46%  "
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.206"
$ws.Range("E42").Value = "  -3.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9994"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.904"
$ws.Range("E44").Value = "  -4.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.76"
$ws.Range("E45").Value = "  -3.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.832"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5997"
$ws.Range("E47").Value = "  -5.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.82"
$ws.Range("E48").Value = "  -5.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.992"
$ws.Range("E49").Value = "  -4.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.159"
$ws.Range("E50").Value = "  -4.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07011"
$ws.Range("E51").Value = "  -6.48%  "
